$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures stored as plain text (e.g. "4.021.53"),
# so force text format before assigning to avoid Excel auto-converting
# them into numeric/date values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.286.91"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.031.31"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.38"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.94"
$ws.Range("E6").Value = "  +2.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.706"
$ws.Range("E7").Value = "  +13.91%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.754"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000326"
$ws.Range("E11").Value = "  -5.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.75"
$ws.Range("E12").Value = "  +11.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.73"
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.662.63"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.023.44"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.13"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.60"
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.990.64"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.10"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "99.08"
$ws.Range("E22").Value = "  +9.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.52"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("E24").Value = "  +4.57%  "
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.09"
$ws.Range("E26").Value = "  -11.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.01"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.63"
$ws.Range("E30").Value = "  +25.98%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "680.56"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.02"
$ws.Range("E34").Value = "  +3.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.49"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "43.07"
$ws.Range("E36").Value = "  +6.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.425"
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0827"
$ws.Range("E39").Value = "  -10.56%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.38"
$ws.Range("E41").Value = "  +7.83%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0488"
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("E45").Value = "  +6.02%  "
$ws.Range("B46").Value = "THORChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.71"
$ws.Range("E46").Value = "  +5.31%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.41"
$ws.Range("E47").Value = "  -3.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.61"
$ws.Range("E48").Value = "  -11.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.01"
$ws.Range("E49").Value = "  -6.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.34"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "145.13"
$ws.Range("E51").Value = "  +1.04%  "
